$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new collection record (row 2): MCH168-1 / STIGEING STEWNFONDS
$ws.Range("A2").Value = "MCH168-1"
$ws.Range("C2").Value = "STIGEING STEWNFONDS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 22C | GRAP COUNT NUMER: NONE"

# Style the new row with the body font (Calibri 10, automatic/theme text color)
$a2 = $ws.Range("A2")
$a2.Font.Name = "Calibri"
$a2.Font.Size = 10
$a2.Font.ThemeColor = 1

# Carry the same formatting across the rest of the row (skip B2 - alternativeIdentifiers
# stays blank/untouched, matching the source record)
$a2.Copy()
$ws.Range("C2:H2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-freeze the header row and select the newly entered record
$ws.Range("A2:J2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
